$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Narrow column J (10th column) from 11 -> 9 (raw OOXML width units)
$ws1.Columns.Item(10).ColumnWidth = 8.17

# Clear a few stray values back to 0
$ws1.Range("D3").Value = 0
$ws1.Range("D4").Value = 0
$ws1.Range("J7").Value = 0
$ws1.Range("N7").Value = 0
$ws1.Range("Q7").Value = 0

# Update the "x de 6" progress counters
$ws1.Range("D8").Value = "0 de 6"
$ws1.Range("J8").Value = "0 de 6"
$ws1.Range("N8").Value = "0 de 6"
$ws1.Range("Q8").Value = "0 de 6"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Column width tweaks
$ws2.Columns.Item(3).ColumnWidth = 11.17
$ws2.Columns.Item(6).ColumnWidth = 10.17

# Roll the month headers forward by one month
$ws2.Range("C1").Value = "abril"
$ws2.Range("D1").Value = "mayo"
$ws2.Range("E1").Value = "junio"
$ws2.Range("F1").Value = "julio"

# Row 3
$ws2.Range("E3").Value = 472.57
$ws2.Range("F3").Value = 0

# Row 4
$ws2.Range("E4").Value = 434.83
$ws2.Range("F4").Value = 0

# Row 5
$ws2.Range("C5").Value = 144.53
$ws2.Range("D5").Value = 11.52
$ws2.Range("E5").Value = 10.44

# Row 6
$ws2.Range("D6").Value = 178.33
$ws2.Range("E6").Value = 0

# Row 7
$ws2.Range("E7").Value = 159.03
$ws2.Range("F7").Value = 0

# Row 8 (totals)
$ws2.Range("C8").Value = 144.53
$ws2.Range("D8").Value = 189.85
$ws2.Range("E8").Value = 1076.87
$ws2.Range("F8").Value = 0
